$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: remove the "Línea de investigación(*)" header (column R)
$ws.Range("R1").Clear()

# Row 2: Area / Subarea1 values change, and "Línea de investigación" value (column R) is removed
$ws.Range("M2").Value = "Ciencias Sociales"
$ws.Range("N2").Value = "ciencias sociales"
$ws.Range("R2").Clear()

# Row 3: Area / Subarea1 / Subarea2 values change, Subarea3 (P) and "Línea de investigación" (R) are removed
$ws.Range("M3").Value = "biociencias"
$ws.Range("N3").Value = "medicina"
$ws.Range("O3").Value = "farmacologia"
$ws.Range("P3").Clear()
$ws.Range("R3").Clear()

# Keep the cursor/selection on the new last data cell (matches reported selection after edit)
$ws.Range("O3").Select() | Out-Null

# Header/footer font tweak ("Normal" -> "Regular")
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Página &P'

